# Update gh-pages to output generated at 456a3b4
# This updates the "想去人数" (want-to-go count) figures in column F
# across the 展览 (sheet1), 演出 (sheet2), 本地生活 (sheet3) and
# 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2590
$ws1.Range("F5").Value  = 914
$ws1.Range("F7").Value  = 1880
$ws1.Range("F8").Value  = 1734
$ws1.Range("F9").Value  = 194
$ws1.Range("F11").Value = 2404
$ws1.Range("F12").Value = 510
$ws1.Range("F13").Value = 187
$ws1.Range("F14").Value = 52
$ws1.Range("F17").Value = 101
$ws1.Range("F18").Value = 8785
$ws1.Range("F20").Value = 6847
$ws1.Range("F21").Value = 11113
$ws1.Range("F23").Value = 190
$ws1.Range("F25").Value = 305
$ws1.Range("F27").Value = 2444
$ws1.Range("F29").Value = 182
$ws1.Range("F30").Value = 2292
$ws1.Range("F31").Value = 435
$ws1.Range("F32").Value = 32
$ws1.Range("F34").Value = 642
$ws1.Range("F35").Value = 284
$ws1.Range("F36").Value = 21
$ws1.Range("F37").Value = 463

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 95
$ws2.Range("F21").Value = 2

# --- 本地生活 sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 88

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 88
$ws4.Range("F7").Value  = 2590
$ws4.Range("F9").Value  = 914
$ws4.Range("F11").Value = 1880
$ws4.Range("F13").Value = 1734
$ws4.Range("F15").Value = 194
$ws4.Range("F17").Value = 510
$ws4.Range("F18").Value = 187
$ws4.Range("F19").Value = 52
$ws4.Range("F22").Value = 101
$ws4.Range("F23").Value = 8785
$ws4.Range("F25").Value = 6847
$ws4.Range("F26").Value = 11113
$ws4.Range("F29").Value = 190
$ws4.Range("F31").Value = 305
$ws4.Range("F37").Value = 32
$ws4.Range("F41").Value = 95
$ws4.Range("F45").Value = 463
$ws4.Range("F47").Value = 2
